$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: extend with two more columns (P1=14, Q1=15), ---
# --- copying the formatting (bold/border/center) from the existing O1 header cell. ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I/K and M/O pairs, and append P/Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
